$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute(", ARQUIVAR E ENCAMINHAR", $true, $false, $false, $false, $false, $true, 1, $false, ", ARQUIVAR E CONFIRMAR DISPONIBILIDADE", 2)
